$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F465").Value = 61723
$ws.Range("G465").Value = 57

$ws.Range("F475").Value = 36720

$ws.Range("F482").Value = 36942
$ws.Range("G482").Value = 24

$ws.Range("F483").Value = 66134

$ws.Range("F514").Value = 7096

$ws.Range("F516").Value = 9455
$ws.Range("F517").Value = 6825
$ws.Range("F518").Value = 7180
$ws.Range("F519").Value = 7979
$ws.Range("F520").Value = 10320

$ws.Range("F523").Value = 10173
$ws.Range("F524").Value = 7811
$ws.Range("F525").Value = 7598
$ws.Range("F526").Value = 8743
$ws.Range("F527").Value = 11385

$ws.Range("F530").Value = 12590
$ws.Range("G530").Value = 41

$ws.Range("F531").Value = 9026

$ws.Range("F532").Value = 10002
$ws.Range("G532").Value = 52

$ws.Range("F533").Value = 11546

$ws.Range("F534").Value = 16229

$ws.Range("F535").Value = 9794
$ws.Range("G535").Value = 22

$ws.Range("F536").Value = 7832
$ws.Range("G536").Value = 39

# New rows 537-539
$ws.Range("A537").Value = 44431
$ws.Range("A537").NumberFormat = "yyyy-mm-dd"
$ws.Range("B537").Value = 394204
$ws.Range("C537").Value = 7578
$ws.Range("D537").Value = 111
$ws.Range("E537").Value = 12547
$ws.Range("F537").Value = 12954
$ws.Range("G537").Value = 50

$ws.Range("A538").Value = 44432
$ws.Range("A538").NumberFormat = "yyyy-mm-dd"
$ws.Range("B538").Value = 394285
$ws.Range("C538").Value = 5380
$ws.Range("D538").Value = 81
$ws.Range("E538").Value = 12547
$ws.Range("F538").Value = 10558
$ws.Range("G538").Value = 28

$ws.Range("A539").Value = 44433
$ws.Range("A539").NumberFormat = "yyyy-mm-dd"
$ws.Range("B539").Value = 394446
$ws.Range("C539").Value = 5931
$ws.Range("D539").Value = 161
$ws.Range("E539").Value = 12547
$ws.Range("F539").Value = 7495
$ws.Range("G539").Value = 31
